$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the table name used in column A (A2:A11) from GroupAddr to tblGroupAddr.
# The D column formulas reference column A and will recalc automatically.
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 -replace [regex]::Escape("[GroupAddr]"), "[tblGroupAddr]"
}

# Correct AddrID values that were wrong.
$ws.Range("C5").Value2 = 9
$ws.Range("C9").Value2 = 8

# Move the active selection to C6 (matches the saved cursor position).
[void]$ws.Range("C6").Select()
